# Auto-generated Excel COM-interop script to update Brynhildr_Profits market data
# Applies per-cell value updates to the ALC, ARM, BSM, CRP, CUL, GSM, LTW sheets
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value2 = 34341.477
$ws.Range("I11").Value2 = 34341.477
$ws.Range("K11").Value2 = 34341.477
$ws.Range("M11").Value2 = -34201.477
$ws.Range("H17").Value2 = 948.5294
$ws.Range("J17").Value2 = 1280.6666
$ws.Range("L17").Value2 = 3841.9998
$ws.Range("N17").Value2 = -4177.9998
$ws.Range("H18").Value2 = 2596.5
$ws.Range("I18").Value2 = 2596.5
$ws.Range("K18").Value2 = 2596.5
$ws.Range("M18").Value2 = -2312.5
$ws.Range("H33").Value2 = 189
$ws.Range("I33").Value2 = 175.125
$ws.Range("K33").Value2 = 175.125
$ws.Range("M33").Value2 = 53.875
$ws.Range("H44").Value2 = 43500
$ws.Range("J44").Value2 = 43500
$ws.Range("L44").Value2 = 43500
$ws.Range("N44").Value2 = -44424
$ws.Range("H62").Value2 = 1851.3334
$ws.Range("I62").Value2 = 1849.75
$ws.Range("J62").Value2 = 1854.5
$ws.Range("K62").Value2 = 1849.75
$ws.Range("L62").Value2 = 1854.5
$ws.Range("M62").Value2 = -1225.75
$ws.Range("N62").Value2 = -3102.5
$ws.Range("H65").Value2 = 1851.3334
$ws.Range("I65").Value2 = 1849.75
$ws.Range("J65").Value2 = 1854.5
$ws.Range("K65").Value2 = 9248.75
$ws.Range("L65").Value2 = 9272.5
$ws.Range("M65").Value2 = -6128.75
$ws.Range("N65").Value2 = -15512.5
$ws.Range("H106").Value2 = 16500
$ws.Range("I106").Value2 = 2400
$ws.Range("J106").Value2 = 18066.666
$ws.Range("K106").Value2 = 2400
$ws.Range("L106").Value2 = 18066.666
$ws.Range("M106").Value2 = -1769
$ws.Range("N106").Value2 = -19328.666
$ws.Range("H135").Value2 = 1209.8276
$ws.Range("I135").Value2 = 1129.8518
$ws.Range("J135").Value2 = 2289.5
$ws.Range("K135").Value2 = 10168.6662
$ws.Range("L135").Value2 = 20605.5
$ws.Range("M135").Value2 = -7633.6662
$ws.Range("N135").Value2 = -25675.5
$ws.Range("H137").Value2 = 17863024
$ws.Range("I137").Value2 = 20835618
$ws.Range("K137").Value2 = 62506854
$ws.Range("M137").Value2 = -62504304
$ws.Range("H141").Value2 = 1053.56
$ws.Range("I141").Value2 = 1053.56
$ws.Range("K141").Value2 = 3160.68
$ws.Range("M141").Value2 = 2019.32

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 641.8333
$ws.Range("I2").Value2 = 695.44446
$ws.Range("J2").Value2 = 481
$ws.Range("K2").Value2 = 695.44446
$ws.Range("L2").Value2 = 481
$ws.Range("M2").Value2 = -582.44446
$ws.Range("N2").Value2 = -707
$ws.Range("H32").Value2 = 697571.6
$ws.Range("I32").Value2 = 773240
$ws.Range("K32").Value2 = 773240
$ws.Range("M32").Value2 = -772953
$ws.Range("H110").Value2 = 880.7222
$ws.Range("I110").Value2 = 880.7222
$ws.Range("K110").Value2 = 880.7222
$ws.Range("M110").Value2 = 1164.2778
$ws.Range("H116").Value2 = 641.8333
$ws.Range("I116").Value2 = 695.44446
$ws.Range("J116").Value2 = 481
$ws.Range("K116").Value2 = 695.44446
$ws.Range("L116").Value2 = 481
$ws.Range("M116").Value2 = 1598.55554
$ws.Range("N116").Value2 = -5069
$ws.Range("H122").Value2 = 2070.25
$ws.Range("I122").Value2 = 1872.3636
$ws.Range("K122").Value2 = 5617.0908
$ws.Range("M122").Value2 = -3167.0908
$ws.Range("H132").Value2 = 2322.476
$ws.Range("I132").Value2 = 1446.5122
$ws.Range("J132").Value2 = 3954.9546
$ws.Range("K132").Value2 = 4339.536599999999
$ws.Range("L132").Value2 = 11864.8638
$ws.Range("M132").Value2 = -1809.536599999999
$ws.Range("N132").Value2 = -16924.8638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 641.8333
$ws.Range("I3").Value2 = 695.44446
$ws.Range("J3").Value2 = 481
$ws.Range("K3").Value2 = 695.44446
$ws.Range("L3").Value2 = 481
$ws.Range("M3").Value2 = -581.44446
$ws.Range("N3").Value2 = -709
$ws.Range("H22").Value2 = 549.6
$ws.Range("I22").Value2 = 549.6
$ws.Range("K22").Value2 = 549.6
$ws.Range("M22").Value2 = -376.6
$ws.Range("H107").Value2 = 1300
$ws.Range("I107").Value2 = 0
$ws.Range("J107").Value2 = 1300
$ws.Range("K107").Value2 = 0
$ws.Range("L107").Value2 = 1300
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value2 = -5140

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1605217.8
$ws.Range("I31").Value2 = 2780844
$ws.Range("K31").Value2 = 2780844
$ws.Range("M31").Value2 = -2780549
$ws.Range("H34").Value2 = 1605217.8
$ws.Range("I34").Value2 = 2780844
$ws.Range("K34").Value2 = 2780844
$ws.Range("M34").Value2 = -2780642
$ws.Range("H58").Value2 = 6637363
$ws.Range("I58").Value2 = 5558166
$ws.Range("J58").Value2 = 14731342
$ws.Range("K58").Value2 = 5558166
$ws.Range("L58").Value2 = 14731342
$ws.Range("M58").Value2 = -5557963
$ws.Range("N58").Value2 = -14731748
$ws.Range("H122").Value2 = 10752.107
$ws.Range("I122").Value2 = 1842.65
$ws.Range("J122").Value2 = 33025.75
$ws.Range("K122").Value2 = 5527.950000000001
$ws.Range("L122").Value2 = 99077.25
$ws.Range("M122").Value2 = -3077.950000000001
$ws.Range("N122").Value2 = -103977.25
$ws.Range("H125").Value2 = 68883.664
$ws.Range("J125").Value2 = 68883.664
$ws.Range("L125").Value2 = 68883.664
$ws.Range("N125").Value2 = -73803.664
$ws.Range("H127").Value2 = 90000
$ws.Range("J127").Value2 = 90000
$ws.Range("L127").Value2 = 90000
$ws.Range("N127").Value2 = -99920
$ws.Range("H136").Value2 = 6637363
$ws.Range("I136").Value2 = 5558166
$ws.Range("J136").Value2 = 14731342
$ws.Range("K136").Value2 = 16674498
$ws.Range("L136").Value2 = 44194026
$ws.Range("M136").Value2 = -16671948
$ws.Range("N136").Value2 = -44199126

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value2 = 315.33334
$ws.Range("I22").Value2 = 273
$ws.Range("K22").Value2 = 819
$ws.Range("M22").Value2 = -650
$ws.Range("H27").Value2 = 315.33334
$ws.Range("I27").Value2 = 273
$ws.Range("K27").Value2 = 819
$ws.Range("M27").Value2 = -717
$ws.Range("H134").Value2 = 1436.2632
$ws.Range("J134").Value2 = 0
$ws.Range("L134").Value2 = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value2 = 3800
$ws.Range("J137").Value2 = 6500
$ws.Range("L137").Value2 = 19500
$ws.Range("N137").Value2 = -29700

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 2025.3889
$ws.Range("I80").Value2 = 1854.3
$ws.Range("K80").Value2 = 1854.3
$ws.Range("M80").Value2 = -856.3
$ws.Range("H83").Value2 = 2025.3889
$ws.Range("I83").Value2 = 1854.3
$ws.Range("K83").Value2 = 9271.5
$ws.Range("M83").Value2 = -4279.5
$ws.Range("H102").Value2 = 2323.9412
$ws.Range("I102").Value2 = 2434.2856
$ws.Range("J102").Value2 = 1809
$ws.Range("K102").Value2 = 2434.2856
$ws.Range("L102").Value2 = 1809
$ws.Range("M102").Value2 = -812.2856000000002
$ws.Range("N102").Value2 = -5053
$ws.Range("H113").Value2 = 2538.625
$ws.Range("I113").Value2 = 2384.8333
$ws.Range("K113").Value2 = 2384.8333
$ws.Range("M113").Value2 = -214.8332999999998
$ws.Range("H122").Value2 = 8599
$ws.Range("I122").Value2 = 11623.25
$ws.Range("J122").Value2 = 4566.6665
$ws.Range("K122").Value2 = 34869.75
$ws.Range("L122").Value2 = 13699.9995
$ws.Range("M122").Value2 = -32419.75
$ws.Range("N122").Value2 = -18599.9995
$ws.Range("H132").Value2 = 16384.809
$ws.Range("I132").Value2 = 18260.479
$ws.Range("K132").Value2 = 54781.437
$ws.Range("M132").Value2 = -52251.437

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4264.2144
$ws.Range("I7").Value2 = 4334.3335
$ws.Range("J7").Value2 = 4138
$ws.Range("K7").Value2 = 4334.3335
$ws.Range("L7").Value2 = 4138
$ws.Range("M7").Value2 = -4222.3335
$ws.Range("N7").Value2 = -4362
$ws.Range("H40").Value2 = 2812.375
$ws.Range("I40").Value2 = 2749.8333
$ws.Range("K40").Value2 = 2749.8333
$ws.Range("M40").Value2 = -2613.8333
$ws.Range("H126").Value2 = 4264.2144
$ws.Range("I126").Value2 = 4334.3335
$ws.Range("J126").Value2 = 4138
$ws.Range("K126").Value2 = 13003.0005
$ws.Range("L126").Value2 = 12414
$ws.Range("M126").Value2 = -10533.0005
$ws.Range("N126").Value2 = -17354
$ws.Range("H133").Value2 = 89297.125
$ws.Range("J133").Value2 = 89297.125
$ws.Range("L133").Value2 = 89297.125
$ws.Range("N133").Value2 = -94357.125
$ws.Range("H136").Value2 = 7145962.5
$ws.Range("I136").Value2 = 4035473.8
$ws.Range("K136").Value2 = 12106421.4
$ws.Range("M136").Value2 = -12103871.4
